$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44964
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 21000
$ws.Range("M2").Value = 20500
$ws.Range("P2").Value = 1139

# Row 3
$ws.Range("D3").Value = 44547
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("P3").Value = 750

# Row 5
$ws.Range("D5").Value = 44984
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17500
$ws.Range("P5").Value = 972

# Row 8
$ws.Range("D8").Value = 44957
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 21000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21500
$ws.Range("P8").Value = 1194

# Row 9
$ws.Range("D9").Value = 44977
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 16500
$ws.Range("L9").Value = 17000
$ws.Range("M9").Value = 16750
$ws.Range("P9").Value = 931
